$d = $word.ActiveDocument

$replacements = @(
    @("675×3=", "705×8="),
    @("219×9=", "524×7="),
    @("189×7=", "936×4="),
    @("427×4=", "392×7="),
    @("625×7=", "440×3="),
    @("313×9=", "944×3="),
    @("117×4=", "744×8="),
    @("439×3=", "534×9="),
    @("596×3=", "201×8="),
    @("645×5=", "608×5="),
    @("946×9=", "930×4="),
    @("875×3=", "618×8="),
    @("657×6=", "638×8="),
    @("727×8=", "556×6="),
    @("302×7=", "966×7="),
    @("390×2=", "715×7="),
    @("372×2=", "343×4="),
    @("627×4=", "878×3="),
    @("152×7=", "824×2="),
    @("884×3=", "713×5="),
    @("935×8=", "963×9="),
    @("395×5=", "802×8="),
    @("581×9=", "360×5="),
    @("242×3=", "189×9="),
    @("815×2=", "711×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
